# ---------------------------------------------------------------------------
# Commit: "Definição dos primeiros requisitos"
#
# 1. Add a new worksheet "Sheet1" at the end of the workbook (after "carro")
#    containing a car-depreciation / price comparison scratchpad.
# 2. Make "RENDA VAARIAVEL" the active sheet/tab again (tabSelected).
# 3. Update the view state of "carro" (scroll position + selection).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Add "Sheet1" after the last existing sheet ("carro") ----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sheet1"

$ws.Range("A1").Value = "valor carro"
$ws.Range("B1").Value = 43000
$ws.Range("D1").Value = "valor carro "
$ws.Range("E1").Value = 70000
$ws.Range("A2").Value = 359
$ws.Range("B2").Formula = "=A2/B1"
$ws.Range("B2").NumberFormat = "0.00%"
$ws.Range("D2").Value = 439
$ws.Range("E2").Formula = "=D2/E1"
$ws.Range("E2").NumberFormat = "0.00%"
$ws.Range("A3").Value = 447
$ws.Range("B3").Formula = "=A3/B1"
$ws.Range("B3").NumberFormat = "0.00%"
$ws.Range("D3").Value = 600
$ws.Range("E3").Formula = "=D3/E1"
$ws.Range("E3").NumberFormat = "0.00%"
$ws.Range("A4").Value = 574
$ws.Range("B4").Formula = "=A4/B1"
$ws.Range("B4").NumberFormat = "0.00%"
$ws.Range("D4").Value = 800
$ws.Range("E4").Formula = "=D4/E1"
$ws.Range("E4").NumberFormat = "0.00%"
$ws.Range("A5").Value = 840
$ws.Range("B5").Formula = "=A5/B1"
$ws.Range("B5").NumberFormat = "0.00%"
$ws.Range("D5").Value = 1135
$ws.Range("E5").Formula = "=D5/E1"
$ws.Range("E5").NumberFormat = "0.00%"
$ws.Range("A6").Value = 482
$ws.Range("B6").Formula = "=A6/B1"
$ws.Range("B6").NumberFormat = "0.00%"
$ws.Range("D6").Value = 671
$ws.Range("E6").Formula = "=D6/E1"
$ws.Range("E6").NumberFormat = "0.00%"
$ws.Range("A7").Value = 632
$ws.Range("B7").Formula = "=A7/B1"
$ws.Range("B7").NumberFormat = "0.00%"
$ws.Range("D7").Value = 955
$ws.Range("E7").Formula = "=D7/E1"
$ws.Range("E7").NumberFormat = "0.00%"
$ws.Range("A8").Value = 482
$ws.Range("B8").Formula = "=A8/B1"
$ws.Range("B8").NumberFormat = "0.00%"
$ws.Range("D8").Value = 671
$ws.Range("E8").Formula = "=D8/E1"
$ws.Range("E8").NumberFormat = "0.00%"
$ws.Range("A9").Value = 840
$ws.Range("B9").Formula = "=A9/B1"
$ws.Range("B9").NumberFormat = "0.00%"
$ws.Range("D9").Value = 1135
$ws.Range("E9").Formula = "=D9/E1"
$ws.Range("E9").NumberFormat = "0.00%"
$ws.Range("A10").Value = 574
$ws.Range("B10").Formula = "=A10/B1"
$ws.Range("B10").NumberFormat = "0.00%"
$ws.Range("D10").Value = 800
$ws.Range("E10").Formula = "=D10/E1"
$ws.Range("E10").NumberFormat = "0.00%"
$ws.Range("A11").Value = 570
$ws.Range("B11").Formula = "=A11/B1"
$ws.Range("B11").NumberFormat = "0.00%"
$ws.Range("D11").Value = 840
$ws.Range("E11").Formula = "=D11/E1"
$ws.Range("E11").NumberFormat = "0.00%"
$ws.Range("A14").Value = "valor do carro"
$ws.Range("B14").Value = "anos "
$ws.Range("C14").Value = "Taxa de depreciação"
$ws.Range("D14").Value = "Valor de depreciação"
$ws.Range("I13").Value = "fordk 1.0 TI - vct"
$ws.Range("J13").Value = 2019
$ws.Range("K13").Value = 43900
$ws.Range("I12").Value = "Carro "
$ws.Range("J12").Value = "Ano"
$ws.Range("K12").Value = "preço"
$ws.Range("L12").Value = "Desvalorização"
$ws.Range("A15").Value = 43000
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 0.2
$ws.Range("C15").NumberFormat = "0.00%"
$ws.Range("D15").Formula = "=A15*C15"
$ws.Range("J14").Value = 2018
$ws.Range("K14").Value = 37852
$ws.Range("L14").Formula = "=1 -K14/K13"
$ws.Range("L14").NumberFormat = "0.00%"
$ws.Range("J15").Value = 2017
$ws.Range("K15").Value = 35554
$ws.Range("L15").Formula = "=1 -K15/K14"
$ws.Range("L15").NumberFormat = "0.00%"
$ws.Range("A16").Formula = "=A15-D15"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 0.1
$ws.Range("C16").NumberFormat = "0.00%"
$ws.Range("D16").Formula = "=A16*C16"
$ws.Range("J16").Value = 2016
$ws.Range("K16").Value = 34231
$ws.Range("L16").Formula = "=1 -K16/K15"
$ws.Range("L16").NumberFormat = "0.00%"
$ws.Range("A17").Formula = "=A16-D16"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 0.1
$ws.Range("C17").NumberFormat = "0.00%"
$ws.Range("D17").Formula = "=A17*C17"
$ws.Range("J17").Value = 2015
$ws.Range("K17").Value = 32503
$ws.Range("L17").Formula = "=1 -K17/K16"
$ws.Range("L17").NumberFormat = "0.00%"
$ws.Range("C18").NumberFormat = "0.00%"
$ws.Range("A21").Formula = "=A15/5"
$ws.Range("B21").Formula = "=A21/5"

# --- 3. Update "carro" sheet view state (scroll + selection) ---------------
$carro = $wb.Worksheets.Item("carro")
$carro.Activate()
$carro.Range("A32").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1

# Give the new "Sheet1" its own selection/view default
$ws.Activate()
$ws.Range("L17").Select()

# --- 2. Restore "RENDA VAARIAVEL" as the active/selected tab ----------------
$wb.Worksheets.Item("RENDA VAARIAVEL").Activate()
